# Adds a "NeedAdjust" flag column (E) to the minutes_format worksheet,
# and highlights the rows that still need manual adjustment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Highlight column B for the flagged rows (except row 159) with a
#    dark red fill. Done first so the new "red fill" style/xf is
#    registered before the new header style below.
# ---------------------------------------------------------------------
$highlightRows = @(3,4,5,6,7,8,20)
foreach ($r in $highlightRows) {
    $ws.Range("B$r").Interior.Color = 192   # RGB(192,0,0) -> FFC00000
}

# ---------------------------------------------------------------------
# 2) New column E: header + values
# ---------------------------------------------------------------------

# Header cell: same look as D1 (bold, centered, top aligned) but with no
# border.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E1").Borders.Item(7).LineStyle = -4142   # xlEdgeLeft   / xlLineStyleNone
$ws.Range("E1").Borders.Item(8).LineStyle = -4142   # xlEdgeTop
$ws.Range("E1").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom
$ws.Range("E1").Borders.Item(10).LineStyle = -4142  # xlEdgeRight
$ws.Range("E1").Value = "NeedAdjust"

# Data cells E2:E230: same look as the rest of the numeric columns
# (style used by column B when its value is 1 - numFmt 0.00, centered,
# light fill).
$ws.Range("B9").Copy() | Out-Null
$ws.Range("E2:E230").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Default everything to 0 ...
$ws.Range("E2:E230").Value = 0

# ... then flag the rows that still need manual adjustment.
$needAdjustRows = @(3,4,5,6,7,8,20,159)
foreach ($r in $needAdjustRows) {
    $ws.Range("E$r").Value = 1
}

# ---------------------------------------------------------------------
# 3) Column E width
# ---------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 14.1

# ---------------------------------------------------------------------
# 4) Sheet view: clear the stale scroll position and select E3:E8
# ---------------------------------------------------------------------
$ws.Range("E3:E8").Select() | Out-Null
